# Auto-generated edit script for grants_solus.xlsx
# Adds 51 new rows (62-112) to the 3rd worksheet ("grants por usuario"):
#   - rows 62-78 : GRANT DELETE, INSERT, SELECT, UPDATE ON TABLE tratamento.tb_pddo_trtmto TO <user>
#   - rows 79-95 : GRANT ALL ON SEQUENCE tratamento.sq_pddo_trtmto TO <user>
#   - rows 96-112: GRANT DELETE, INSERT, SELECT, UPDATE ON TABLE tratamento.tb_c_cid TO <user>
# for the same list of 17 users/roles used throughout the sheet (plus two
# newly introduced ones: "vanessa" and "aline").
#
# The column order in which cells are written below (all of column A for
# the first two new GRANT phrases, then column B for those same rows, then
# column A/B for the third phrase) mirrors how the workbook's shared string
# table ends up populated (tb_pddo_trtmto, sq_pddo_trtmto, vanessa, aline,
# tb_c_cid) so the resulting file matches the authored edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

$names = @(
    "`"adriana.paes`"",
    "`"kleverson.antonio`"",
    "`"luana.mourao`"",
    "`"mariana.brider`"",
    "`"adriana.paes`"",
    "`"larissa.rocha`"",
    "`"livia.hallack`"",
    "`"vanessa.cirilo`"",
    "`"milena.ribeiral`"",
    "`"tatiane.neto`"",
    "`"tatyene.nehrer`"",
    "`"victor.quinet`"",
    "`"lidice.lenz`"",
    "`"alberlania.muller`"",
    "`"marcia.costa`"",
    "vanessa",
    "aline"
)

# --- Column A -------------------------------------------------------
# Category 1: rows 62-78
for ($r = 62; $r -le 78; $r++) {
    $ws.Cells.Item($r, 1).Value = "GRANT DELETE, INSERT, SELECT, UPDATE ON TABLE tratamento.tb_pddo_trtmto TO"
}
# Category 2: rows 79-95
for ($r = 79; $r -le 95; $r++) {
    $ws.Cells.Item($r, 1).Value = "GRANT ALL ON SEQUENCE tratamento.sq_pddo_trtmto TO"
}

# --- Column B (categories 1 & 2) ------------------------------------
for ($i = 0; $i -lt $names.Length; $i++) {
    $ws.Cells.Item(62 + $i, 2).Value = $names[$i]
    $ws.Cells.Item(79 + $i, 2).Value = $names[$i]
}

# --- Column A (category 3): rows 96-112 ------------------------------
for ($r = 96; $r -le 112; $r++) {
    $ws.Cells.Item($r, 1).Value = "GRANT DELETE, INSERT, SELECT, UPDATE ON TABLE tratamento.tb_c_cid TO"
}

# --- Column B (category 3) -------------------------------------------
for ($i = 0; $i -lt $names.Length; $i++) {
    $ws.Cells.Item(96 + $i, 2).Value = $names[$i]
}

# --- Column C: literal ";" for every new row --------------------------
for ($r = 62; $r -le 112; $r++) {
    $ws.Cells.Item($r, 3).Value = ";"
}

# --- Formatting: column B uses the same style as the rest of the sheet
# (Arial 9pt, wrap text, vertically centered) -- copy it from an existing
# formatted cell so no new style entries are introduced.
$ws.Range("B61").Copy() | Out-Null
$ws.Range("B62:B112").PasteSpecial(-4122) | Out-Null

# --- Column D formulas -------------------------------------------------
# Rows 62-65 extend the original shared formula that covered D2:D61.
for ($r = 62; $r -le 65; $r++) {
    $ws.Cells.Item($r, 4).Formula = "=A$r&`" `"&B$r&`" `"&C$r"
}
# Rows 66-112 form a new shared formula block.
$ws.Range("D66:D112").Formula = "=A66&`" `"&B66&`" `"&C66"

# --- View state: select the tail of the newly added data, matching the
# cursor position left behind by the edit.
$ws.Activate()
$ws.Range("D95:D112").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 88
$excel.ActiveWindow.ScrollColumn = 2

$wb.Save()
